$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Increase search space size values
$ws.Range("B2").Value = 1000
$ws.Range("D2").Value = 100
$ws.Range("E2").Value = 100

$ws.Range("B3").Value = 5000
$ws.Range("D3").Value = 500
$ws.Range("E3").Value = 100

$ws.Range("B4").Value = 10000
$ws.Range("D4").Value = 1000
$ws.Range("E4").Value = 100

# Update the active cell selection
$ws.Range("C12").Select()
